# Auto-generated edit script for ctrl-q-master-items.xlsx
# Implements: "feat: Include color info when importing master dimensions from Excel file"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1): add Color / Per value color columns ---
$ws.Range("G1").Value = "Color"
$ws.Range("H1").Value = "Per value color"

# --- Update existing data rows 2-8 (measures): tidy label/expression text ---
# Row 2
$ws.Range("A2").Value = "measure"
$ws.Range("B2").Value = "No. of sold units"
$ws.Range("C2").Value = "Number of units sold during selected time period."
$ws.Range("D2").Value = "'='Sold units'"
$ws.Range("E2").Value = "Sales"
$ws.Range("F2").Value = "'=Sum(UnitsInOrder)"

# Row 3
$ws.Range("A3").Value = "measure"
$ws.Range("B3").Value = "No. of sold units (LY)"
$ws.Range("C3").Value = "Number of units sold last year."
$ws.Range("D3").Value = "'='Sold units LY'"
$ws.Range("E3").Value = "Sales, LY"
$ws.Range("F3").Value = "Sum(UnitsInOrder_LY)"

# Row 4
$ws.Range("A4").Value = "measure"
$ws.Range("B4").Value = "Revenue EUR"
$ws.Range("C4").Value = "Revenue during selected time period."
$ws.Range("D4").Value = "'='Revenue'"
$ws.Range("E4").Value = "Sales"
$ws.Range("F4").Value = "Sum(Revenue)"

# Row 5
$ws.Range("A5").Value = "measure"
$ws.Range("B5").Value = "Revenue EUR (LY)"
$ws.Range("C5").Value = "Revenue during last year."
$ws.Range("D5").Value = "'='Revenue LY'"
$ws.Range("E5").Value = "Sales, LY"
$ws.Range("F5").Value = "Sum(Revenue_LY)"

# Row 6
$ws.Range("A6").Value = "measure"
$ws.Range("B6").Value = "Profit EUR"
$ws.Range("C6").Value = "Profit during selected time period."
$ws.Range("D6").Value = "'='Profit'"
$ws.Range("E6").Value = "Sales"
$ws.Range("F6").Value = "Sum(Profit)"

# Row 7
$ws.Range("A7").Value = "measur"
$ws.Range("B7").Value = "Profit USD"
$ws.Range("C7").Value = "Profit during selected time period."
$ws.Range("D7").Value = "'='Profit'"
$ws.Range("E7").Value = "Sales"
$ws.Range("F7").Value = "Sum(Profit)"

# Row 8
$ws.Range("A8").Value = "measure"
$ws.Range("B8").Value = "Profit EUR (LY)"
$ws.Range("C8").Value = "Profit during last year"
$ws.Range("D8").Value = "'='Profit LY'"
$ws.Range("E8").Value = "Sales, LY"
$ws.Range("F8").Value = "Sum(Profit_LY)"

# --- Row 9 (Country dimension): update text + add color JSON in G/H ---
$ws.Range("A9").Value = "dim-single"
$ws.Range("B9").Value = "Country"
$ws.Range("C9").Value = "Country where a unit was sold."
$ws.Range("D9").Value = "'='Country sold'"
$ws.Range("E9").Value = "Geo, DimCat1"
$ws.Range("F9").Value = "Country"
$ws.Range("B9").Copy()
$ws.Range("G9:H9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G9").Value = "{`n  `"baseColor`": {`n    `"color`": `"#bbbbbb`",`n    `"index`": -1`n  }`n}"
$ws.Range("H9").Value = "{`n    `"colors`": [`n        {`n            `"value`": `"Afghanistan`",`n            `"baseColor`": {`n                `"color`": `"#8a85c6`",`n                `"index`": -1`n            }`n        },`n        {`n            `"value`": `"Albania`",`n            `"baseColor`": {`n                `"color`": `"#aaaaaa`",`n                `"index`": -1`n            }`n        },`n        {`n            `"value`": `"Algeria`",`n            `"baseColor`": {`n                `"color`": `"#a16090`",`n                `"index`": 9`n            }`n        }`n    ],`n    `"nul`": {`n        `"color`": `"#c8c7a9`",`n        `"index`": 16`n    },`n    `"oth`": {`n        `"color`": `"#ffec6e`",`n        `"index`": -1`n    },`n    `"pal`": null,`n    `"single`": null,`n    `"usePal`": true,`n    `"autoFill`": true`n}"

# --- Row 10 (Sales month dimension): update text, clear Tag, add color JSON in G ---
$ws.Range("A10").Value = "dim-single"
$ws.Range("B10").Value = "Sales month"
$ws.Range("C10").Value = "Date in which a unit was sold."
$ws.Range("D10").Value = "'='Sales month'"
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = "Month_Sales"
$ws.Range("B10").Copy()
$ws.Range("G10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G10").Value = "{`n  `"baseColor`": {`n    `"color`": `"#bbbbbb`",`n    `"index`": -1`n  }`n}"

# --- Row 11 (Salesperson dimension): update text, add color JSON in H ---
$ws.Range("A11").Value = "dim-single"
$ws.Range("B11").Value = "Salesperson"
$ws.Range("C11").Value = "The person who sold the unit."
$ws.Range("D11").Value = "'='Salesperson'"
$ws.Range("E11").Value = "Staff, Sales"
$ws.Range("F11").Value = "Salesperson"
$ws.Range("B11").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H11").Value = "{`n    `"colors`": [`n        {`n            `"value`": `"Afghanistan`",`n            `"baseColor`": {`n                `"color`": `"#8a85c6`",`n                `"index`": -1`n            }`n        },`n        {`n            `"value`": `"Albania`",`n            `"baseColor`": {`n                `"color`": `"#aaaaaa`",`n                `"index`": -1`n            }`n        },`n        {`n            `"value`": `"Algeria`",`n            `"baseColor`": {`n                `"color`": `"#a16090`",`n                `"index`": 9`n            }`n        }`n    ],`n    `"nul`": {`n        `"color`": `"#c8c7a9`",`n        `"index`": 16`n    },`n    `"oth`": {`n        `"color`": `"#ffec6e`",`n        `"index`": -1`n    },`n    `"pal`": null,`n    `"single`": null,`n    `"usePal`": true,`n    `"autoFill`": true`n}"

# --- Row 12 (new Color dimension): copy formatting from row 11, then set values ---
$ws.Range("A11:F11").Copy()
$ws.Range("A12:F12").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A12").Value = "dim-single"
$ws.Range("B12").Value = "Color"
$ws.Range("C12").Value = "Color of sold unit"
$ws.Range("D12").Value = "'='Unit color'"
$ws.Range("E12").Value = "Sales, Color"
$ws.Range("F12").Value = "UnitColor"

# --- Row heights ---
$ws.Rows.Item(9).RowHeight = 159
$ws.Rows.Item(10).RowHeight = 90
$ws.Rows.Item(11).RowHeight = 409.6

# --- Column widths (ColumnWidth chars are offset by 5/6 from stored OOXML width) ---
$ws.Columns.Item(1).ColumnWidth = 15.498697916666666
$ws.Columns.Item(3).ColumnWidth = 30.666666666666668
$ws.Columns.Item(6).ColumnWidth = 19.998697916666668
$ws.Columns.Item(7).ColumnWidth = 17.830729166666668
$ws.Columns.Item(8).ColumnWidth = 22.330729166666668

# --- Selection / view ---
$ws.Range("A12").Select()
$excel.ActiveWindow.ScrollRow = 11
